$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ------------------------------------------------------------------
# 1) Insert a new "Compile-fail tests for a Call." list item right
#    after "Make better use of auto.", moving the _GoBack bookmark
#    (which trailed the old last item) onto the new last item.
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq "Make better use of auto.`r") {
        $target = $cand
        break
    }
}

$r = $d.Range($target.Range.Start, $target.Range.End)
$xml1 = "<w:p xmlns:w='$wNs' w:rsidR='009720EE' w:rsidRDefault='009720EE' w:rsidP='00144F10'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>Make better use of auto.</w:t></w:r></w:p>" + `
         "<w:p xmlns:w='$wNs'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>Compile-fail tests for a Call.</w:t></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>"
$r.InsertXML($xml1)

# ------------------------------------------------------------------
# 2) Move the rendered-page-break marker: it now falls before "VEH
#    hooking (both INT3 and DR)." instead of before "Transactional
#    hooking." (a new earlier paragraph pushed the page break back).
# ------------------------------------------------------------------
$pVeh = $null
$pTrans = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $t = $cand.Range.Text
    if ($t -eq "VEH hooking (both INT3 and DR).`r") { $pVeh = $cand }
    if ($t -eq "Transactional hooking.`r") { $pTrans = $cand }
}

$r2 = $d.Range($pVeh.Range.Start, $pTrans.Range.End)
$xml2 = "<w:p xmlns:w='$wNs' w:rsidR='00312CFE' w:rsidRPr='008B06FC' w:rsidRDefault='00312CFE' w:rsidP='00312CFE'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r w:rsidRPr='008B06FC'><w:lastRenderedPageBreak/><w:t>VEH hooking (both INT3 and DR).</w:t></w:r></w:p>" + `
         "<w:p xmlns:w='$wNs' w:rsidR='00312CFE' w:rsidRPr='008B06FC' w:rsidRDefault='00312CFE' w:rsidP='00312CFE'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r w:rsidRPr='008B06FC'><w:t>Transactional hooking.</w:t></w:r></w:p>"
$r2.InsertXML($xml2)
